$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "<four>"
$ws.Range("C2").Value = 53

$ws.Range("B3").Value = "<time>"
$ws.Range("C3").Value = 56

$ws.Range("B4").Value = "<by>"
$ws.Range("C4").Value = 49

$ws.Range("C5").Value = 54

$ws.Range("B6").Value = "<enter>"
$ws.Range("C6").Value = 49

$ws.Range("C7").Value = 47

$ws.Range("B8").Value = "<be>"
$ws.Range("C8").Value = 51

$ws.Range("C9").Value = 47

$ws.Range("B10").Value = "<them>"
$ws.Range("C10").Value = 46

$ws.Range("B11").Value = "<on>"
$ws.Range("C11").Value = 51

$ws.Range("B12").Value = "<them>"
$ws.Range("C12").Value = 52

$ws.Range("C13").Value = 51

$ws.Range("C14").Value = 54

$ws.Range("B16").Value = "<lima>"
$ws.Range("C16").Value = 49

$ws.Range("B17").Value = "<controw>"
$ws.Range("C17").Value = 60

$ws.Range("C18").Value = 45
